# This edit inserts one new weekly price-report data row for
# "Macroferia Regional de Talca" (Acelga) ahead of the existing row 116,
# shifting all subsequent rows (old 116-167) down by one (new 117-168).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 116; everything from the old row 116
# onward (including row 167) shifts down by one row, to 117..168.
$ws.Rows.Item(116).Insert()

# Populate the newly inserted row 116 with the new weekly record.
$ws.Cells.Item(116, 1).Value  = 5
$ws.Cells.Item(116, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(116, 3).Value  = "Maule"
$ws.Cells.Item(116, 4).Value  = 44455
$ws.Cells.Item(116, 5).Value  = 7
$ws.Cells.Item(116, 6).Value  = 100112009
$ws.Cells.Item(116, 7).Value  = "Acelga"
$ws.Cells.Item(116, 8).Value  = "Sin especificar"
$ws.Cells.Item(116, 9).Value  = "Primera"
$ws.Cells.Item(116, 10).Value = 500
$ws.Cells.Item(116, 11).Value = 2500
$ws.Cells.Item(116, 12).Value = 2500
$ws.Cells.Item(116, 13).Value = 2500
$ws.Cells.Item(116, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(116, 15).Value = "Región del Maule"
$ws.Cells.Item(116, 16).Value = 625
$ws.Cells.Item(116, 17).Value = 4
$ws.Cells.Item(116, 18).Value = "Hortaliza"
